$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2,4) '80.297.53'
$ws.Cells.Item(2,5).Value = '  +5.09%  '

# Row 3
Set-TextValue $ws.Cells.Item(3,4) '3.209.13'
$ws.Cells.Item(3,5).Value = '  +4.01%  '

# Row 4
$ws.Cells.Item(4,5).Value = '  -0.03%  '

# Row 5
Set-TextValue $ws.Cells.Item(5,4) '210.10'
$ws.Cells.Item(5,5).Value = '  +5.85%  '

# Row 6
Set-TextValue $ws.Cells.Item(6,4) '632.75'
$ws.Cells.Item(6,5).Value = '  +2.68%  '

# Row 7
$ws.Cells.Item(7,5).Value = '  +32.45%  '

# Row 8
Set-TextValue $ws.Cells.Item(8,4) '0.999'
$ws.Cells.Item(8,5).Value = '  -0.08%  '

# Row 9
Set-TextValue $ws.Cells.Item(9,4) '0.604'
$ws.Cells.Item(9,5).Value = '  +9.75%  '

# Row 10
Set-TextValue $ws.Cells.Item(10,4) '3.208.55'
$ws.Cells.Item(10,5).Value = '  +3.95%  '

# Row 11
Set-TextValue $ws.Cells.Item(11,4) '0.619'
$ws.Cells.Item(11,5).Value = '  +40.93%  '

# Row 12
Set-TextValue $ws.Cells.Item(12,4) '0.0000262'
$ws.Cells.Item(12,5).Value = '  +35.32%  '

# Row 13
$ws.Cells.Item(13,5).Value = '  +3.20%  '

# Row 14
Set-TextValue $ws.Cells.Item(14,4) '5.40'
$ws.Cells.Item(14,5).Value = '  +3.46%  '

# Row 15
Set-TextValue $ws.Cells.Item(15,4) '3.803.29'
$ws.Cells.Item(15,5).Value = '  +5.36%  '

# Row 16
Set-TextValue $ws.Cells.Item(16,4) '32.47'
$ws.Cells.Item(16,5).Value = '  +11.90%  '

# Row 17
Set-TextValue $ws.Cells.Item(17,4) '80.367.73'
$ws.Cells.Item(17,5).Value = '  +5.16%  '

# Row 18
Set-TextValue $ws.Cells.Item(18,4) '3.217.55'
$ws.Cells.Item(18,5).Value = '  +4.94%  '

# Row 19
Set-TextValue $ws.Cells.Item(19,4) '14.57'
$ws.Cells.Item(19,5).Value = '  +6.97%  '

# Row 20
Set-TextValue $ws.Cells.Item(20,4) '447.89'
$ws.Cells.Item(20,5).Value = '  +17.43%  '

# Row 21
Set-TextValue $ws.Cells.Item(21,4) '9.33'
$ws.Cells.Item(21,5).Value = '  +1.96%  '

# Row 22
$ws.Cells.Item(22,5).Value = '  +21.51%  '

# Row 23
Set-TextValue $ws.Cells.Item(23,4) '5.35'
$ws.Cells.Item(23,5).Value = '  +21.54%  '

# Row 24
Set-TextValue $ws.Cells.Item(24,4) '6.81'
$ws.Cells.Item(24,5).Value = '  +5.04%  '

# Row 25
Set-TextValue $ws.Cells.Item(25,4) '3.385.81'
$ws.Cells.Item(25,5).Value = '  +4.53%  '

# Row 26
Set-TextValue $ws.Cells.Item(26,4) '77.57'
$ws.Cells.Item(26,5).Value = '  +7.20%  '

# Row 27
$ws.Cells.Item(27,5).Value = '  +10.27%  '

# Row 28
Set-TextValue $ws.Cells.Item(28,4) '10.95'
$ws.Cells.Item(28,5).Value = '  +10.83%  '

# Row 29
$ws.Cells.Item(29,5).Value = '  -0.25%  '

# Row 30
$ws.Cells.Item(30,5).Value = '  +15.42%  '

# Row 31
Set-TextValue $ws.Cells.Item(31,4) '9.27'
$ws.Cells.Item(31,5).Value = '  +11.75%  '

# Row 32
Set-TextValue $ws.Cells.Item(32,4) '0.999'
$ws.Cells.Item(32,5).Value = '  -0.42%  '

# Row 33
Set-TextValue $ws.Cells.Item(33,4) '556.75'
$ws.Cells.Item(33,5).Value = '  +10.85%  '

# Row 34
Set-TextValue $ws.Cells.Item(34,4) '1.50'
$ws.Cells.Item(34,5).Value = '  +6.57%  '

# Row 35
Set-TextValue $ws.Cells.Item(35,4) '0.153'
$ws.Cells.Item(35,5).Value = '  +23.18%  '

# Row 36
Set-TextValue $ws.Cells.Item(36,4) '2.03'
$ws.Cells.Item(36,5).Value = '  +6.11%  '

# Row 37
Set-TextValue $ws.Cells.Item(37,4) '23.70'
$ws.Cells.Item(37,5).Value = '  +14.43%  '

# Row 38
$ws.Cells.Item(38,5).Value = '  +21.78%  '

# Row 39
Set-TextValue $ws.Cells.Item(39,4) '0.421'
$ws.Cells.Item(39,5).Value = '  +11.17%  '

# Row 40
Set-TextValue $ws.Cells.Item(40,4) '1.00'
$ws.Cells.Item(40,5).Value = '  +0.03%  '

# Row 41
Set-TextValue $ws.Cells.Item(41,4) '165.76'
$ws.Cells.Item(41,5).Value = '  +1.83%  '

# Row 42: full replace (coin swap)
$ws.Cells.Item(42,2).Value = 'WhiteBITCoin'
$ws.Cells.Item(42,3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Cells.Item(42,4) '20.62'
$ws.Cells.Item(42,5).Value = '  +2.80%  '

# Row 43: full replace (coin swap)
$ws.Cells.Item(43,2).Value = 'RenderToken'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Cells.Item(43,4) '5.75'
$ws.Cells.Item(43,5).Value = '  +12.36%  '

# Row 44
Set-TextValue $ws.Cells.Item(44,4) '192.13'
$ws.Cells.Item(44,5).Value = '  -0.55%  '

# Row 45
$ws.Cells.Item(45,5).Value = '  +0.02%  '

# Row 46
$ws.Cells.Item(46,5).Value = '  +11.71%  '

# Row 47
$ws.Cells.Item(47,5).Value = '  +12.41%  '

# Row 48
Set-TextValue $ws.Cells.Item(48,4) '0.797'
$ws.Cells.Item(48,5).Value = '  +0.55%  '

# Row 49
$ws.Cells.Item(49,5).Value = '  +7.40%  '

# Row 50
Set-TextValue $ws.Cells.Item(50,4) '43.58'
$ws.Cells.Item(50,5).Value = '  +5.77%  '

# Row 51
Set-TextValue $ws.Cells.Item(51,4) '4.33'
$ws.Cells.Item(51,5).Value = '  +11.53%  '
